$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from E1 (header style) to F1 so F1 matches existing header formatting
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F1").Value = "time_taken"
$ws.Range("F2").Value = "2021-10-05 13:38:53.227928"
$ws.Range("F3").Value = "2021-10-05 13:38:53.227938"
$ws.Range("F4").Value = "2021-10-05 13:38:53.227941"
$ws.Range("F5").Value = "2021-10-05 13:38:53.227943"
$ws.Range("F6").Value = "2021-10-05 13:38:53.227945"
$ws.Range("F7").Value = "2021-10-05 13:38:53.227947"
$ws.Range("F8").Value = "2021-10-05 13:38:53.227949"
$ws.Range("F9").Value = "2021-10-05 13:38:53.227951"
$ws.Range("F10").Value = "2021-10-05 13:38:53.227954"
$ws.Range("F11").Value = "2021-10-05 13:38:53.228001"
$ws.Range("F12").Value = "2021-10-05 13:38:53.228007"
$ws.Range("F13").Value = "2021-10-05 13:38:53.228010"
$ws.Range("F14").Value = "2021-10-05 13:38:53.228013"
$ws.Range("F15").Value = "2021-10-05 13:38:53.228017"
$ws.Range("F16").Value = "2021-10-05 13:38:53.228021"
$ws.Range("F17").Value = "2021-10-05 13:38:53.228024"
$ws.Range("F18").Value = "2021-10-05 13:38:53.228028"
$ws.Range("F19").Value = "2021-10-05 13:38:53.228058"
$ws.Range("F20").Value = "2021-10-05 13:38:53.228061"
$ws.Range("F21").Value = "2021-10-05 13:38:53.228082"
$ws.Range("F22").Value = "2021-10-05 13:38:53.228084"
$ws.Range("F23").Value = "2021-10-05 13:38:53.228086"
$ws.Range("F24").Value = "2021-10-05 13:38:53.228088"
$ws.Range("F25").Value = "2021-10-05 13:38:53.228090"
$ws.Range("F26").Value = "2021-10-05 13:38:53.228093"
$ws.Range("F27").Value = "2021-10-05 13:38:53.228095"
$ws.Range("F28").Value = "2021-10-05 13:38:53.228097"
$ws.Range("F29").Value = "2021-10-05 13:38:53.228099"
$ws.Range("F30").Value = "2021-10-05 13:38:53.228101"
$ws.Range("F31").Value = "2021-10-05 13:38:53.228104"
$ws.Range("F32").Value = "2021-10-05 13:38:53.228106"
$ws.Range("F33").Value = "2021-10-05 13:38:53.228108"
$ws.Range("F34").Value = "2021-10-05 13:38:53.228110"
$ws.Range("F35").Value = "2021-10-05 13:38:53.228113"
$ws.Range("F36").Value = "2021-10-05 13:38:53.228115"
$ws.Range("F37").Value = "2021-10-05 13:38:53.228117"
$ws.Range("F38").Value = "2021-10-05 13:38:53.228119"
$ws.Range("F39").Value = "2021-10-05 13:38:53.228121"
$ws.Range("F40").Value = "2021-10-05 13:38:53.228123"
$ws.Range("F41").Value = "2021-10-05 13:38:53.228125"
$ws.Range("F42").Value = "2021-10-05 13:38:53.228127"
$ws.Range("F43").Value = "2021-10-05 13:38:53.228129"
$ws.Range("F44").Value = "2021-10-05 13:38:53.228131"
$ws.Range("F45").Value = "2021-10-05 13:38:53.228133"
$ws.Range("F46").Value = "2021-10-05 13:38:53.228135"
$ws.Range("F47").Value = "2021-10-05 13:38:53.228137"
$ws.Range("F48").Value = "2021-10-05 13:38:53.228139"
$ws.Range("F49").Value = "2021-10-05 13:38:53.228141"
$ws.Range("F50").Value = "2021-10-05 13:38:53.228143"
$ws.Range("F51").Value = "2021-10-05 13:38:53.228145"
$ws.Range("F52").Value = "2021-10-05 13:38:53.228147"
$ws.Range("F53").Value = "2021-10-05 13:38:53.228148"
$ws.Range("F54").Value = "2021-10-05 13:38:53.228151"
$ws.Range("F55").Value = "2021-10-05 13:38:53.228153"
$ws.Range("F56").Value = "2021-10-05 13:38:53.228155"
$ws.Range("F57").Value = "2021-10-05 13:38:53.228156"
$ws.Range("F58").Value = "2021-10-05 13:38:53.228159"
$ws.Range("F59").Value = "2021-10-05 13:38:53.228161"
$ws.Range("F60").Value = "2021-10-05 13:38:53.228162"
$ws.Range("F61").Value = "2021-10-05 13:38:53.228164"
$ws.Range("F62").Value = "2021-10-05 13:38:53.228166"
$ws.Range("F63").Value = "2021-10-05 13:38:53.228168"
$ws.Range("F64").Value = "2021-10-05 13:38:53.228170"
$ws.Range("F65").Value = "2021-10-05 13:38:53.228172"
$ws.Range("F66").Value = "2021-10-05 13:38:53.228175"
$ws.Range("F67").Value = "2021-10-05 13:38:53.228177"
$ws.Range("F68").Value = "2021-10-05 13:38:53.228179"
$ws.Range("F69").Value = "2021-10-05 13:38:53.228181"
$ws.Range("F70").Value = "2021-10-05 13:38:53.228183"
$ws.Range("F71").Value = "2021-10-05 13:38:53.228185"
$ws.Range("F72").Value = "2021-10-05 13:38:53.228187"
$ws.Range("F73").Value = "2021-10-05 13:38:53.228189"
$ws.Range("F74").Value = "2021-10-05 13:38:53.228191"
$ws.Range("F75").Value = "2021-10-05 13:38:53.228193"
$ws.Range("F76").Value = "2021-10-05 13:38:53.228195"
$ws.Range("F77").Value = "2021-10-05 13:38:53.228197"
$ws.Range("F78").Value = "2021-10-05 13:38:53.228201"
$ws.Range("F79").Value = "2021-10-05 13:38:53.228203"
$ws.Range("F80").Value = "2021-10-05 13:38:53.228205"
$ws.Range("F81").Value = "2021-10-05 13:38:53.228207"
$ws.Range("F82").Value = "2021-10-05 13:38:53.228209"
$ws.Range("F83").Value = "2021-10-05 13:38:53.228211"
$ws.Range("F84").Value = "2021-10-05 13:38:53.228213"
$ws.Range("F85").Value = "2021-10-05 13:38:53.228215"
$ws.Range("F86").Value = "2021-10-05 13:38:53.228217"
$ws.Range("F87").Value = "2021-10-05 13:38:53.228219"
$ws.Range("F88").Value = "2021-10-05 13:38:53.228221"
$ws.Range("F89").Value = "2021-10-05 13:38:53.228223"
$ws.Range("F90").Value = "2021-10-05 13:38:53.228226"
$ws.Range("F91").Value = "2021-10-05 13:38:53.228228"
$ws.Range("F92").Value = "2021-10-05 13:38:53.228230"
$ws.Range("F93").Value = "2021-10-05 13:38:53.228232"
$ws.Range("F94").Value = "2021-10-05 13:38:53.228235"
$ws.Range("F95").Value = "2021-10-05 13:38:53.228238"
$ws.Range("F96").Value = "2021-10-05 13:38:53.228240"
$ws.Range("F97").Value = "2021-10-05 13:38:53.228242"
$ws.Range("F98").Value = "2021-10-05 13:38:53.228244"
